# Weekly update: insert two new daily price records for
# "Feria Lagunitas de Puerto Montt - Frutilla" right after the existing
# row for date 44567 (old row 229), pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 230:231 - this shifts every row from the old
# 230 onward down by two (old 230 -> new 232, ..., old 309 -> new 311),
# matching the new dimension A1:T311.
$ws.Rows("230:231").Insert()

# --- New row 230 -------------------------------------------------------
$ws.Range("A230").Value = 4
$ws.Range("B230").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C230").Value = "Los Lagos"
$ws.Range("D230").Value = 44900
$ws.Range("E230").Value = 10
$ws.Range("F230").Value = "Fruta"
$ws.Range("G230").Value = 100101
$ws.Range("H230").Value = "Berries"
$ws.Range("I230").Value = 100112025
$ws.Range("J230").Value = "Frutilla"
$ws.Range("K230").Value = "Sin especificar"
$ws.Range("L230").Value = "Primera"
$ws.Range("M230").Value = 400
$ws.Range("N230").Value = 10000
$ws.Range("O230").Value = 11000
$ws.Range("P230").Value = 10500
$ws.Range("Q230").Value = "$/bandeja 7 kilos"
$ws.Range("R230").Value = "Provincia de Melipilla"
$ws.Range("S230").Value = 1500
$ws.Range("T230").Value = 7

# --- New row 231 -------------------------------------------------------
$ws.Range("A231").Value = 4
$ws.Range("B231").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C231").Value = "Los Lagos"
$ws.Range("D231").Value = 44900
$ws.Range("E231").Value = 10
$ws.Range("F231").Value = "Fruta"
$ws.Range("G231").Value = 100101
$ws.Range("H231").Value = "Berries"
$ws.Range("I231").Value = 100112025
$ws.Range("J231").Value = "Frutilla"
$ws.Range("K231").Value = "Sin especificar"
$ws.Range("L231").Value = "Primera"
$ws.Range("M231").Value = 400
$ws.Range("N231").Value = 10000
$ws.Range("O231").Value = 11000
$ws.Range("P231").Value = 10500
$ws.Range("Q231").Value = "$/caja 7 kilos"
$ws.Range("R231").Value = "Región de La Araucanía"
$ws.Range("S231").Value = 1500
$ws.Range("T231").Value = 7
